$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Punto medio Circulo")
$ws.Activate()
